$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions) - update "想去人数" (interested count) column F
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 6433
$ws1.Range("F5").Value = 383
$ws1.Range("F9").Value = 87
$ws1.Range("F15").Value = 3140
$ws1.Range("F18").Value = 1813

# Sheet "全部类型" (all types) - same underlying rows, mirrored with an extra row
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 6433
$ws4.Range("F5").Value = 383
$ws4.Range("F10").Value = 87
$ws4.Range("F16").Value = 3140
$ws4.Range("F19").Value = 1813
